# StringExceptionBug project re-save.
# The only semantic change in this commit is the text typed into cell B18
# on Sheet1: it was "&gt;" and is now "@" (shared-string table grows as a
# side effect of Excel recording the new string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B18").Value = "@"
